$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1 (table row index 1)
$t.Cell(1,1).Range.Text = "99÷2=49, 1"
$t.Cell(1,2).Range.Text = "29÷7=4, 1"
$t.Cell(1,3).Range.Text = "97÷2=48, 1"
$t.Cell(1,4).Range.Text = "23÷4=5, 3"
$t.Cell(1,5).Range.Text = "27÷3=9, 0"

# Row 5 (table row index 5)
$t.Cell(5,1).Range.Text = "34÷6=5, 4"
$t.Cell(5,2).Range.Text = "84÷6=14, 0"
$t.Cell(5,3).Range.Text = "40÷8=5, 0"
$t.Cell(5,4).Range.Text = "80÷4=20, 0"
$t.Cell(5,5).Range.Text = "69÷3=23, 0"

# Row 9 (table row index 9)
$t.Cell(9,1).Range.Text = "76÷6=12, 4"
$t.Cell(9,2).Range.Text = "23÷9=2, 5"
$t.Cell(9,3).Range.Text = "82÷4=20, 2"
$t.Cell(9,4).Range.Text = "63÷5=12, 3"
$t.Cell(9,5).Range.Text = "27÷7=3, 6"

# Row 13 (table row index 13)
$t.Cell(13,1).Range.Text = "46÷4=11, 2"
$t.Cell(13,2).Range.Text = "51÷3=17, 0"
$t.Cell(13,3).Range.Text = "87÷5=17, 2"
$t.Cell(13,4).Range.Text = "90÷4=22, 2"
$t.Cell(13,5).Range.Text = "99÷4=24, 3"

# Row 17 (table row index 17)
$t.Cell(17,1).Range.Text = "14÷9=1, 5"
$t.Cell(17,2).Range.Text = "59÷4=14, 3"
$t.Cell(17,3).Range.Text = "86÷7=12, 2"
$t.Cell(17,4).Range.Text = "30÷2=15, 0"
$t.Cell(17,5).Range.Text = "10÷2=5, 0"
